$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = 956
$ws.Range("C3").Value = 2007
$ws.Range("D3").Value = 3771
$ws.Range("E3").Value = 8551
$ws.Range("F3").Value = 11300
$ws.Range("G3").Value = 34500

$ws.Range("B4").Value = 125.82912
$ws.Range("C4").Value = 263.192576
$ws.Range("D4").Value = 493.879296
$ws.Range("E4").Value = 1120.927744
$ws.Range("F4").Value = 1477.443584
$ws.Range("G4").Value = 4516.216832

$ws.Range("B5").Value = 1044.61
$ws.Range("C5").Value = 974.4400000000001
$ws.Range("D5").Value = 1041.58
$ws.Range("E5").Value = 846.63
$ws.Range("F5").Value = 1359.72
$ws.Range("G5").Value = 915.13

$ws.Range("B6").Value = 1844
$ws.Range("C6").Value = 1778
$ws.Range("D6").Value = 1860
$ws.Range("E6").Value = 1991
$ws.Range("F6").Value = 2933
$ws.Range("G6").Value = 1729

$ws.Range("B7").Value = 3195
$ws.Range("C7").Value = 1909
$ws.Range("D7").Value = 2040
$ws.Range("E7").Value = 2278
$ws.Range("F7").Value = 4047
$ws.Range("G7").Value = 2311

$ws.Range("B12").Value = 4460
$ws.Range("C12").Value = 7964
$ws.Range("D12").Value = 16800
$ws.Range("E12").Value = 30800
$ws.Range("F12").Value = 66300
$ws.Range("G12").Value = 85800

$ws.Range("B13").Value = 18.2452224
$ws.Range("C13").Value = 32.6107136
$ws.Range("D13").Value = 68.5768704
$ws.Range("E13").Value = 125.82912
$ws.Range("F13").Value = 271.581184
$ws.Range("G13").Value = 351.27296

$ws.Range("B14").Value = 223.97017
$ws.Range("C14").Value = 247.32561
$ws.Range("D14").Value = 237.47531
$ws.Range("E14").Value = 258.67722
$ws.Range("F14").Value = 240.26036
$ws.Range("G14").Value = 371.73891

$ws.Range("B15").Value = 1564.672
$ws.Range("C15").Value = 1499.136
$ws.Range("D15").Value = 1482.752
$ws.Range("E15").Value = 1515.52
$ws.Range("F15").Value = 1548.288
$ws.Range("G15").Value = 1810.432

$ws.Range("B16").Value = 1679.36
$ws.Range("C16").Value = 1662.976
$ws.Range("D16").Value = 1613.824
$ws.Range("E16").Value = 1679.36
$ws.Range("F16").Value = 1744.896
$ws.Range("G16").Value = 3293.184

$ws.Range("B21").Value = 7937
$ws.Range("C21").Value = 10900
$ws.Range("D21").Value = 16200
$ws.Range("E21").Value = 17100
$ws.Range("F21").Value = 17600
$ws.Range("G21").Value = 18100

$ws.Range("B22").Value = 1040.187392
$ws.Range("C22").Value = 1428.160512
$ws.Range("D22").Value = 2122.317824
$ws.Range("E22").Value = 2246.049792
$ws.Range("F22").Value = 2312.11008
$ws.Range("G22").Value = 2378.170368

$ws.Range("B23").Value = 66.55
$ws.Range("C23").Value = 64.92
$ws.Range("D23").Value = 69.94
$ws.Range("E23").Value = 94.31999999999999
$ws.Range("F23").Value = 163.22
$ws.Range("G23").Value = 343.23

$ws.Range("B24").Value = 112
$ws.Range("C24").Value = 95
$ws.Range("D24").Value = 90
$ws.Range("E24").Value = 123
$ws.Range("F24").Value = 225
$ws.Range("G24").Value = 486

$ws.Range("B25").Value = 235
$ws.Range("C25").Value = 208
$ws.Range("D25").Value = 155
$ws.Range("E25").Value = 217
$ws.Range("F25").Value = 330
$ws.Range("G25").Value = 1037

$ws.Range("B30").Value = 106000
$ws.Range("C30").Value = 197000
$ws.Range("D30").Value = 311000
$ws.Range("E30").Value = 406000
$ws.Range("F30").Value = 401000
$ws.Range("G30").Value = 437000

$ws.Range("B31").Value = 436.207616
$ws.Range("C31").Value = 808.452096
$ws.Range("D31").Value = 1275.068416
$ws.Range("E31").Value = 1665.138688
$ws.Range("F31").Value = 1642.070016
$ws.Range("G31").Value = 1788.870656

$ws.Range("B32").Value = $null
$ws.Range("C32").Value = $null
$ws.Range("D32").Value = $null
$ws.Range("E32").Value = 10.36
$ws.Range("F32").Value = 16.04
$ws.Range("G32").Value = 26.14

$ws.Range("B33").Value = 10
$ws.Range("C33").Value = 10
$ws.Range("D33").Value = 11
$ws.Range("E33").Value = 14.4
$ws.Range("F33").Value = 32
$ws.Range("G33").Value = 41

$ws.Range("B34").Value = 11
$ws.Range("C34").Value = 11
$ws.Range("D34").Value = 12
$ws.Range("E34").Value = 16.32
$ws.Range("F34").Value = 52
$ws.Range("G34").Value = 75

$ws.Range("B39").Value = 20500
$ws.Range("C39").Value = 28800
$ws.Range("D39").Value = 43100
$ws.Range("E39").Value = 42000
$ws.Range("F39").Value = 38600
$ws.Range("G39").Value = 36900

$ws.Range("B40").Value = 2684.35456
$ws.Range("C40").Value = 3781.165056
$ws.Range("D40").Value = 5650.776064
$ws.Range("E40").Value = 5506.072576
$ws.Range("F40").Value = 5064.62208
$ws.Range("G40").Value = 4830.789632

$ws.Range("B41").Value = 47.97
$ws.Range("C41").Value = 66.77
$ws.Range("D41").Value = 90.79000000000001
$ws.Range("E41").Value = 183.59
$ws.Range("F41").Value = 403.84
$ws.Range("G41").Value = 853.66

$ws.Range("B42").Value = 94
$ws.Range("C42").Value = 151
$ws.Range("D42").Value = 306
$ws.Range("E42").Value = 529
$ws.Range("F42").Value = 947
$ws.Range("G42").Value = 1926

$ws.Range("B43").Value = 215
$ws.Range("C43").Value = 297
$ws.Range("D43").Value = 506
$ws.Range("E43").Value = 627
$ws.Range("F43").Value = 1057
$ws.Range("G43").Value = 2671

$ws.Range("B48").Value = 655000
$ws.Range("C48").Value = 819000
$ws.Range("D48").Value = 1049000
$ws.Range("E48").Value = 1036000
$ws.Range("F48").Value = 1051000
$ws.Range("G48").Value = 1074000

$ws.Range("B49").Value = 2684.35456
$ws.Range("C49").Value = 3355.4432
$ws.Range("D49").Value = 4294.967296
$ws.Range("E49").Value = 4243.587072
$ws.Range("F49").Value = 4303.355904
$ws.Range("G49").Value = 4400.873472

$ws.Range("B50").Value = 1.41504
$ws.Range("C50").Value = 2.03072
$ws.Range("D50").Value = 3.5334
$ws.Range("E50").Value = 7.40739
$ws.Range("F50").Value = 14.11403
$ws.Range("G50").Value = 28.67039

$ws.Range("B51").Value = 0.49
$ws.Range("C51").Value = 0.47
$ws.Range("D51").Value = 0.498
$ws.Range("E51").Value = 0.51
$ws.Range("F51").Value = 0.524
$ws.Range("G51").Value = 0.628

$ws.Range("B52").Value = 61.184
$ws.Range("C52").Value = 64.768
$ws.Range("D52").Value = 160.768
$ws.Range("E52").Value = 432.128
$ws.Range("F52").Value = 815.104
$ws.Range("G52").Value = 1630.208

$ws.Range("B57").Value = 6481
$ws.Range("C57").Value = 9183
$ws.Range("D57").Value = 11200
$ws.Range("E57").Value = 12700
$ws.Range("F57").Value = 12300
$ws.Range("G57").Value = 13100

$ws.Range("B58").Value = 849.34656
$ws.Range("C58").Value = 1203.765248
$ws.Range("D58").Value = 1471.152128
$ws.Range("E58").Value = 1659.895808
$ws.Range("F58").Value = 1615.855616
$ws.Range("G58").Value = 1712.324608

$ws.Range("B59").Value = 66.89
$ws.Range("C59").Value = 62.45
$ws.Range("D59").Value = 69.42
$ws.Range("E59").Value = 92.59999999999999
$ws.Range("F59").Value = 164.1
$ws.Range("G59").Value = 343.27

$ws.Range("B60").Value = 114
$ws.Range("C60").Value = 75
$ws.Range("D60").Value = 82
$ws.Range("E60").Value = 115
$ws.Range("F60").Value = 210
$ws.Range("G60").Value = 457

$ws.Range("B61").Value = 141
$ws.Range("C61").Value = 108
$ws.Range("D61").Value = 90
$ws.Range("E61").Value = 121
$ws.Range("F61").Value = 258
$ws.Range("G61").Value = 627

$ws.Range("B66").Value = 96900
$ws.Range("C66").Value = 178000
$ws.Range("D66").Value = 226000
$ws.Range("E66").Value = 295000
$ws.Range("F66").Value = 329000
$ws.Range("G66").Value = 359000

$ws.Range("B67").Value = 397.410304
$ws.Range("C67").Value = 729.808896
$ws.Range("D67").Value = 925.892608
$ws.Range("E67").Value = 1206.910976
$ws.Range("F67").Value = 1345.323008
$ws.Range("G67").Value = 1470.103552

$ws.Range("B68").Value = $null
$ws.Range("C68").Value = $null
$ws.Range("D68").Value = $null
$ws.Range("E68").Value = $null
$ws.Range("F68").Value = 15.23
$ws.Range("G68").Value = 25.12

$ws.Range("B69").Value = 8.768000000000001
$ws.Range("C69").Value = 10
$ws.Range("D69").Value = 11
$ws.Range("E69").Value = 13.12
$ws.Range("F69").Value = 28.8
$ws.Range("G69").Value = 37

$ws.Range("B70").Value = 9.407999999999999
$ws.Range("C70").Value = 10
$ws.Range("D70").Value = 12
$ws.Range("E70").Value = 14.784
$ws.Range("F70").Value = 45.312
$ws.Range("G70").Value = 63

Write-Host "Applied rfuse results updates"